# Generate Report for Handback
# Update the "Correspond Handoff Datetime" and "Correspond Handback DateTime"
# timestamps on the zh-cn and de-de report sheets.

$wb = $excel.ActiveWorkbook

# --- zh-cn sheet ---
$zhcn = $wb.Worksheets.Item("zh-cn")
# Correspond Handoff Datetime (column D), rows 2 and 3
$zhcn.Range("D2").Value = "2016-02-22 18:09:22"
$zhcn.Range("D3").Value = "2016-02-22 18:09:22"
# Correspond Handback DateTime (column G), rows 2 and 3
$zhcn.Range("G2").Value = "2016-02-22 18:10:12"
$zhcn.Range("G3").Value = "2016-02-22 18:10:12"

# --- de-de sheet ---
$dede = $wb.Worksheets.Item("de-de")
# Correspond Handoff Datetime (column D), rows 2 and 3
$dede.Range("D2").Value = "2016-02-22 18:09:33"
$dede.Range("D3").Value = "2016-02-22 18:09:33"
# Correspond Handback DateTime (column G), rows 2 and 3
$dede.Range("G2").Value = "2016-02-22 18:10:37"
$dede.Range("G3").Value = "2016-02-22 18:10:37"
